$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 19 (the "zero.system.document" / 文档管理 menu row) entirely.
# This shifts rows 20:35 up to 19:34 and drops the now-unused shared
# strings that only that row referenced.
$ws.Rows("19").Delete() | Out-Null

# Re-apply the two shared formulas that straddled the deleted row so they
# stay grouped as shared formulas (Excel keeps them shared after a normal
# row-delete; re-asserting the formula on the post-shift range reproduces
# that behaviour exactly).
$ws.Range("B21:B22").Formula = "=A$20"
$ws.Range("B33:B34").Formula = "=A$30"

# Excel leaves the newly-shifted row 19 selected (as an entire row) after
# deleting the row at that position.
$ws.Rows("19").Select() | Out-Null
